$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Hora" (G) values advance from 3 to 4 for rows 2-51
$ws.Range("G2:G51").Value = "'4"
$ws.Range("G2:G51").Style = "Normal"

# Updated "Price" (D) values for the rows that changed
$ws.Range("D2").Value = "'266.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.258"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'3.566"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.535"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.416"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8249"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1645"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08250"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03541"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03190"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09186"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.767"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001635"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04683"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.006418"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006186"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.0001502"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.725"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'2.256"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.01372"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3288"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.04684"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007020"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.004306"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1117"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01055"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006229"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0009906"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.9910"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Value = "'0.00001903"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01242"
$ws.Range("D51").Style = "Normal"
